$wb = $excel.ActiveWorkbook

# --- Sheet: Resumen_Calidad ---
$ws1 = $wb.Worksheets.Item("Resumen_Calidad")

$ws1.Range("B2").Value = 34

$ws1.Range("C2").Value = "{'N°': 0, 'Actividad del Proyecto': 26, '¿A qué actor va dirigida?': 26, 'Número de Beneficiarios': 26, 'Entrega Dotación (SI / NO)': 26, 'Descripción de la Dotación Entregada': 32, 'Evidencia de la Actividad': 34, 'Evidencia_URL': 34, 'Observaciones Generales': 0, 'Hoja': 0, 'Nombre_Proyecto': 26, 'Nombre Proyecto': 0, 'FUENTES': 0, 'PROYECTOS': 0, 'BENEFICIARIOS': 0, 'Unnamed: 4': 0, 'Unnamed: 5': 0}"

$ws1.Range("C4").Value = "{'N°': 0, 'Actividad del Proyecto': 0, 'Total Ejecutado': 0, 'Componente PAM': 63, '¿A qué actor va dirigida?': 66, 'Número de Beneficiarios': 76, 'Entrega Dotación (SI / NO)': 70, 'Descripción de la Dotación Entregada': 84, 'Evidencia de la Actividad': 0, 'Evidencia_URL': 0, 'Observaciones Generales': 0, 'Hoja': 0, 'Nombre_Proyecto': 0, 'Nombre Proyecto': 0, 'FUENTES': 0, 'PROYECTOS': 0, 'BENEFICIARIOS': 0, 'Unnamed: 4': 0, 'Unnamed: 5': 0}"

$ws1.Range("C6").Value = "{'N°': 0, 'DANE IEO': 0, 'MUNICIO': 0, 'NOMBRE_IEO': 0, '# Directivos Beneficiados': 219, '# Administrativos Beneficiados': 209, '# Docentes Beneficiados': 228, '# Estudiantes Beneficiados': 234, '¿Recibió Asistencia Técnica?': 203, 'Modalidad de la Asistencia Técnica': 234, '¿Recibió Dotación?': 225, 'Dotación Recibida': 228, 'Hoja': 0, 'Nombre Proyecto': 0, 'FUENTES': 0, 'PROYECTOS': 0, 'BENEFICIARIOS': 0, 'Unnamed: 4': 0, 'Unnamed: 5': 0}"

$ws1.Range("C8").Value = "{'N°': 0, 'DANE IEO': 0, 'MUNICIO': 0, 'NOMBRE_IEO': 0, '# Directivos Beneficiados': 1605, '# Administrativos Beneficiados': 1691, '# Docentes Beneficiados': 1577, '# Estudiantes Beneficiados': 1870, '¿Recibió Asistencia Técnica?': 1614, 'Modalidad de la Asistencia Técnica': 1614, '¿Recibió Dotación?': 1612, 'Dotación Recibida': 1777, 'Hoja': 0, '# Padres - Madres y Cuidadores': 2056, 'Nombre Proyecto': 0, 'FUENTES': 0, 'PROYECTOS': 0, 'BENEFICIARIOS': 0, 'Unnamed: 4': 0, 'Unnamed: 5': 0}"

$ws1.Range("C12").Value = "{'Col1': 14, 'N°': 0, 'Nombre Proyecto': 0, 'Código BPIN': 0, 'Código PI': 0, 'Apropiación Definitiva 2025': 0, 'Adición': 0, 'Total Ejecutado 2025': 0, 'Difrencia \nApro - Ejec': 0, 'Porcentaje de Ejecución': 0, 'RECURSOS': 0, 'Responsable': 2, 'Enlace Técnico': 5, 'ENLACE PROYECTO': 14, 'ENLACE BENEFICIARIOS': 14, 'Documentos del Proyecto': 0, 'IGP': 14, 'Avance en el Cargue de información': 0, 'Col19': 14, 'Hoja': 0}"

# --- Sheet: Resumen_Depuracion ---
$ws3 = $wb.Worksheets.Item("Resumen_Depuracion")
$ws3.Range("B2").Value = 34
$ws3.Range("C2").Value = 34
